$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.669.43'
$ws.Range('E2').Value = '  +1.93%  '
$ws.Range('D3').Value = '3.924.19'
$ws.Range('E3').Value = '  +0.30%  '
$ws.Range('D5').Value = '''533.33'
$ws.Range('E5').Value = '  +9.79%  '
$ws.Range('D6').Value = '''144.49'
$ws.Range('E6').Value = '  -0.80%  '
$ws.Range('E7').Value = '  -0.99%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  -0.63%  '
$ws.Range('D10').Value = '''0.173'
$ws.Range('E10').Value = '  +4.21%  '
$ws.Range('D11').Value = '''0.0000334'
$ws.Range('E11').Value = '  -3.31%  '
$ws.Range('D12').Value = '''42.46'
$ws.Range('E12').Value = '  -1.45%  '
$ws.Range('D13').Value = '4.556.64'
$ws.Range('E13').Value = '  +0.42%  '
$ws.Range('D14').Value = '''10.30'
$ws.Range('E14').Value = '  -3.94%  '
$ws.Range('D15').Value = '3.926.41'
$ws.Range('E15').Value = '  -0.13%  '
$ws.Range('E16').Value = '  +8.62%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').Value = '''0.136'
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').Value = '''13.91'
$ws.Range('E18').Value = '  -2.97%  '
$ws.Range('D19').Value = '''19.86'
$ws.Range('E19').Value = '  -0.61%  '
$ws.Range('D20').Value = '69.542.99'
$ws.Range('E20').Value = '  +1.67%  '
$ws.Range('D21').Value = '''439.35'
$ws.Range('E21').Value = '  +1.63%  '
$ws.Range('D22').Value = '''3.36'
$ws.Range('E22').Value = '  -3.42%  '
$ws.Range('D23').Value = '''14.43'
$ws.Range('E23').Value = '  -4.86%  '
$ws.Range('D24').Value = '''4.15'
$ws.Range('E24').Value = '  +13.57%  '
$ws.Range('D25').Value = '''88.16'
$ws.Range('E25').Value = '  +0.13%  '
$ws.Range('D26').Value = '''11.59'
$ws.Range('E26').Value = '  +0.98%  '
$ws.Range('D27').Value = '''10.73'
$ws.Range('E27').Value = '  -4.71%  '
$ws.Range('D28').Value = '''36.52'
$ws.Range('E28').Value = '  -3.51%  '
$ws.Range('D29').Value = '''697.66'
$ws.Range('E29').Value = '  -2.73%  '
$ws.Range('D30').Value = '''13.23'
$ws.Range('E30').Value = '  -4.12%  '
$ws.Range('D31').Value = '''0.126'
$ws.Range('E31').Value = '  -2.64%  '
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').Value = '''2.84'
$ws.Range('E32').Value = '  -3.43%  '
$ws.Range('B33').Value = 'OKB'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D33').Value = '''68.69'
$ws.Range('E33').Value = '  +12.64%  '
$ws.Range('D34').Value = '''0.449'
$ws.Range('E34').Value = '  +14.52%  '
$ws.Range('E35').Value = '  -2.58%  '
$ws.Range('D36').Value = '''40.33'
$ws.Range('E36').Value = '  -2.49%  '
$ws.Range('D37').Value = '0.0₃0843'
$ws.Range('E37').Value = '  -3.69%  '
$ws.Range('E38').Value = '  +2.31%  '
$ws.Range('D39').Value = '''0.999'
$ws.Range('E39').Value = '  -0.05%  '
$ws.Range('E40').Value = '  -0.09%  '
$ws.Range('D41').Value = '''0.0484'
$ws.Range('E41').Value = '  -3.61%  '
$ws.Range('E42').Value = '  +3.13%  '
$ws.Range('D43').Value = '''2.76'
$ws.Range('E43').Value = '  -8.69%  '
$ws.Range('D44').Value = '''2.96'
$ws.Range('E44').Value = '  -4.77%  '
$ws.Range('D45').Value = '''3.15'
$ws.Range('E45').Value = '  +11.65%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').Value = '''3.39'
$ws.Range('E46').Value = '  +0.12%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').Value = '''0.142'
$ws.Range('E47').Value = '  -0.20%  '
$ws.Range('D48').Value = '''3.32'
$ws.Range('E48').Value = '  -2.73%  '
$ws.Range('D49').Value = '0.0₆0342'
$ws.Range('E49').Value = '  +1.60%  '
$ws.Range('D50').Value = '''145.01'
$ws.Range('D51').Value = '''2.07'
$ws.Range('E51').Value = '  -2.79%  '
